# Revert "RESTORE: Recover all 973 original multi-industry template files"
# for the IT_KPI_Dashboard.xlsx template: put the Information Technology
# wording back in place of the Artificial Intelligence / Machine Learning
# wording that the restore commit had reintroduced.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Instructions & User Guide"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")

$ws1.Range("A1").Value = "Information Technology KPI Dashboard - User Guide & Instructions"

$ws1.Range("A19").Value = "System Uptime Percentage"
$ws1.Range("B19").Value = "Key performance indicator for cloud computing and infrastructure automation"

$ws1.Range("B22").Value = "Key performance indicator for cloud computing and infrastructure automation"

$ws1.Range("B24").Value = "Key performance indicator for cloud computing and infrastructure automation"

# ---------------------------------------------------------------------
# Sheet 2: "KPI Dashboard"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("KPI Dashboard")

$ws2.Range("A1").Value = "INFORMATION TECHNOLOGY - KPI DASHBOARD"
$ws2.Range("A2").Value = "Project: Cloud Infrastructure Migration"

# KPI001 (row 8)
$ws2.Range("B8").Value = "System Uptime Percentage"
$ws2.Range("I8").Value = "Chief Technology Officer"
$ws2.Range("K8").Value = "Critical KPI for Information Technology success"

# KPI002 (row 9)
$ws2.Range("B9").Value = "Mean Time To Repair (MTTR)"
$ws2.Range("I9").Value = "IT Managers"
$ws2.Range("K9").Value = "Critical KPI for Information Technology success"

# KPI003 (row 10)
$ws2.Range("B10").Value = "Deployment Frequency"
$ws2.Range("I10").Value = "DevOps Engineers"
$ws2.Range("K10").Value = "Critical KPI for Information Technology success"

# KPI004 (row 11)
$ws2.Range("I11").Value = "System Administrators"
$ws2.Range("K11").Value = "Critical KPI for Information Technology success"

# KPI005 (row 12)
$ws2.Range("K12").Value = "Critical KPI for Information Technology success"

# KPI006 (row 13)
$ws2.Range("K13").Value = "Critical KPI for Information Technology success"

# KPI007 (row 14)
$ws2.Range("I14").Value = "Chief Technology Officer"
$ws2.Range("K14").Value = "Critical KPI for Information Technology success"

# KPI008 (row 15)
$ws2.Range("I15").Value = "IT Managers"
$ws2.Range("K15").Value = "Critical KPI for Information Technology success"

# KPI009 (row 16)
$ws2.Range("I16").Value = "DevOps Engineers"
$ws2.Range("K16").Value = "Critical KPI for Information Technology success"

# KPI010 (row 17)
$ws2.Range("I17").Value = "System Administrators"
$ws2.Range("K17").Value = "Critical KPI for Information Technology success"

# KPI011 (row 18)
$ws2.Range("K18").Value = "Critical KPI for Information Technology success"

# KPI012 (row 19)
$ws2.Range("K19").Value = "Critical KPI for Information Technology success"

# KPI013 (row 20)
$ws2.Range("I20").Value = "Chief Technology Officer"
$ws2.Range("K20").Value = "Critical KPI for Information Technology success"

# KPI014 (row 21)
$ws2.Range("I21").Value = "IT Managers"
$ws2.Range("K21").Value = "Critical KPI for Information Technology success"

# KPI015 (row 22)
$ws2.Range("I22").Value = "DevOps Engineers"
$ws2.Range("K22").Value = "Critical KPI for Information Technology success"
